$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Add new row of data (row 11): STT, Ten, SL, Gia -> E11 already has shared formula C11*D11
$ws.Range("A11").Value = 10
$ws.Range("B11").Value = "Bánh bi"
$ws.Range("C11").Value = 1
$ws.Range("D11").Value = 55500

# Update the selected cell to H9, matching the saved selection in the diff
$ws.Range("H9").Select()

$wb.Application.Calculate()
